{"js": "// Replace the date line and each \"AxB=C\" multiplication-table answer\n// with its updated value, following strict document order so that\n// duplicate-looking numbers (e.g. two different cells that happen to\n// reduce to the same digits) never get cross-matched.\nconst replacements = [\n  [\"2025-02-07 Friday\", \"2025-02-08 Saturday\"],\n  [\"74\u00d740=2960\", \"62\u00d744=2728\"],\n  [\"64\u00d734=2176\", \"75\u00d736=2700\"],\n  [\"81\u00d785=6885\", \"28\u00d746=1288\"],\n  [\"67\u00d712=804\", \"75\u00d755=4125\"],\n  [\"21\u00d729=609\", \"52\u00d787=4524\"],\n  [\"68\u00d763=4284\", \"29\u00d760=1740\"],\n  [\"13\u00d747=611\", \"37\u00d791=3367\"],\n  [\"29\u00d768=1972\", \"64\u00d796=6144\"],\n  [\"35\u00d713=455\", \"21\u00d724=504\"],\n  [\"71\u00d773=5183\", \"22\u00d748=1056\"],\n  [\"16\u00d783=1328\", \"92\u00d764=5888\"],\n  [\"28\u00d741=1148\", \"15\u00d719=285\"],\n  [\"49\u00d749=2401\", \"88\u00d717=1496\"],\n  [\"61\u00d795=5795\", \"69\u00d748=3312\"],\n  [\"30\u00d771=2130\", \"61\u00d798=5978\"],\n  [\"34\u00d788=2992\", \"90\u00d798=8820\"],\n  [\"51\u00d745=2295\", \"47\u00d792=4324\"],\n  [\"59\u00d737=2183\", \"50\u00d725=1250\"],\n  [\"95\u00d721=1995\", \"35\u00d771=2485\"],\n  [\"37\u00d727=999\", \"73\u00d721=1533\"],\n  [\"91\u00d715=1365\", \"95\u00d738=3610\"],\n  [\"78\u00d714=1092\", \"54\u00d782=4428\"],\n  [\"17\u00d757=969\", \"25\u00d788=2200\"],\n  [\"91\u00d796=8736\", \"52\u00d788=4576\"],\n  [\"29\u00d721=609\", \"93\u00d732=2976\"],\n];\n\nconst body = context.document.body;\n\nfor (const [before, after] of replacements) {\n  const found = body.search(before, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  if (found.items.length === 0) {\n    throw new Error(`Could not find text to replace: \"${before}\"`);\n  }\n\n  // Only the first occurrence is expected (all source strings are unique),\n  // but guard against accidental repeats by only touching the first hit.\n  found.items[0].insertText(after, \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# Update the date line and each two-digit-by-two-digit multiplication\n# \"AxB=C\" answer to its new value. Replacements are applied in strict\n# document order (date first, then each table cell, row by row, left to\n# right) and each one only touches its single, exact match (MatchCase\n# on, MatchWholeWord off since the needle already is the whole cell\n# text, Replace = wdReplaceOne) so that no other cell can accidentally\n# be touched.\n$d = $word.ActiveDocument\n\n$replacements = @(\n  @(\"2025-02-07 Friday\", \"2025-02-08 Saturday\"),\n  @(\"74\u00d740=2960\", \"62\u00d744=2728\"),\n  @(\"64\u00d734=2176\", \"75\u00d736=2700\"),\n  @(\"81\u00d785=6885\", \"28\u00d746=1288\"),\n  @(\"67\u00d712=804\", \"75\u00d755=4125\"),\n  @(\"21\u00d729=609\", \"52\u00d787=4524\"),\n  @(\"68\u00d763=4284\", \"29\u00d760=1740\"),\n  @(\"13\u00d747=611\", \"37\u00d791=3367\"),\n  @(\"29\u00d768=1972\", \"64\u00d796=6144\"),\n  @(\"35\u00d713=455\", \"21\u00d724=504\"),\n  @(\"71\u00d773=5183\", \"22\u00d748=1056\"),\n  @(\"16\u00d783=1328\", \"92\u00d764=5888\"),\n  @(\"28\u00d741=1148\", \"15\u00d719=285\"),\n  @(\"49\u00d749=2401\", \"88\u00d717=1496\"),\n  @(\"61\u00d795=5795\", \"69\u00d748=3312\"),\n  @(\"30\u00d771=2130\", \"61\u00d798=5978\"),\n  @(\"34\u00d788=2992\", \"90\u00d798=8820\"),\n  @(\"51\u00d745=2295\", \"47\u00d792=4324\"),\n  @(\"59\u00d737=2183\", \"50\u00d725=1250\"),\n  @(\"95\u00d721=1995\", \"35\u00d771=2485\"),\n  @(\"37\u00d727=999\", \"73\u00d721=1533\"),\n  @(\"91\u00d715=1365\", \"95\u00d738=3610\"),\n  @(\"78\u00d714=1092\", \"54\u00d782=4428\"),\n  @(\"17\u00d757=969\", \"25\u00d788=2200\"),\n  @(\"91\u00d796=8736\", \"52\u00d788=4576\"),\n  @(\"29\u00d721=609\", \"93\u00d732=2976\"),\n)\n\nforeach ($pair in $replacements) {\n  $findText = $pair[0]\n  $replaceText = $pair[1]\n\n  $rng = $d.Content\n  $found = $rng.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n\n  if (-not $found) {\n    throw \"Could not find text to replace: $findText\"\n  }\n}\n"}
